# [ADDITIONAL SCRAPING] added code to scrape more data about a player's
# batting performance in a match, also updated the excel sheets.
#
# 1. Insert a new "Player Info" worksheet before the existing "ODI Batting"
#    sheet, with player bio fields (ID, NAME, BATTING_HAND, BOWL_STYLE).
# 2. Rename the "ODI Batting" sheet's MATCH_CARD_LINK column to MATCH_CODE
#    and replace the full scorecard URLs with the bare numeric match code.

$wb = $excel.ActiveWorkbook

# --- 1. Add the new "Player Info" sheet ---
# Worksheets.Add() with no args inserts the new sheet immediately before
# the currently active sheet, i.e. right in front of "ODI Batting" (the
# only/active sheet in this workbook).
$infoSheet = $wb.Worksheets.Add()
$infoSheet.Name = "Player Info"

# Header row
$infoSheet.Range("A1").Value = "ID"
$infoSheet.Range("B1").Value = "NAME"
$infoSheet.Range("C1").Value = "BATTING_HAND"
$infoSheet.Range("D1").Value = "BOWL_STYLE"

# Match the bold / bordered / centered-top header style used on the
# "ODI Batting" sheet's own header row.
$infoHeader = $infoSheet.Range("A1:D1")
$infoHeader.Font.Bold = $true
$infoHeader.Borders.LineStyle = 1
$infoHeader.HorizontalAlignment = -4108
$infoHeader.VerticalAlignment = -4160

# Data row - force text so the numeric-looking ID stays a string, matching
# how every other data cell in this workbook is stored.
$infoSheet.Range("A2").Value = "'6469"
$infoSheet.Range("B2").Value = "Riaz Hassan"
$infoSheet.Range("C2").Value = "Right Handed"
$infoSheet.Range("D2").Value = "Does Not Bowl | Unknown"

# --- 2. Update the "ODI Batting" sheet's MATCH_CARD_LINK column ---
# Re-fetch by name now, after the insert, so the reference points at the
# right sheet object (sheet references resolve by current position).
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$battingSheet.Range("D1").Value = "MATCH_CODE"

# Leading apostrophe forces text storage (like the rest of the sheet's
# inline-string cells) instead of letting Excel coerce these into numbers.
$battingSheet.Range("D2").Value = "'4530"
$battingSheet.Range("D3").Value = "'4538"
$battingSheet.Range("D4").Value = "'4539"
